$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 74, shifting existing rows 74-106 down to 75-107
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with the new data record
$ws.Cells.Item(74, 1).Value2 = 6
$ws.Cells.Item(74, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(74, 3).Value2 = "Metropolitana"
$ws.Cells.Item(74, 4).Value2 = 44845
$ws.Cells.Item(74, 5).Value2 = 13
$ws.Cells.Item(74, 6).Value2 = 100114007
$ws.Cells.Item(74, 7).Value2 = "Jengibre"
$ws.Cells.Item(74, 8).Value2 = "Sin especificar"
$ws.Cells.Item(74, 9).Value2 = "Primera"
$ws.Cells.Item(74, 10).Value2 = 800
$ws.Cells.Item(74, 11).Value2 = 11000
$ws.Cells.Item(74, 12).Value2 = 12000
$ws.Cells.Item(74, 13).Value2 = 11562
$ws.Cells.Item(74, 14).Value2 = '$/caja 13 kilos'
$ws.Cells.Item(74, 15).Value2 = "Perú"
$ws.Cells.Item(74, 16).Value2 = 889
$ws.Cells.Item(74, 17).Value2 = 13
$ws.Cells.Item(74, 18).Value2 = "Hortaliza"
